# The "Dados" sheet (first sheet, xl/worksheets/sheet1.xml) has column A
# manually widened - originally 10.85546875 "OOXML characters" wide, the
# author dragged it out to 16.85546875 wide so the longer values in that
# sheet fit after the newest Powell/CDEEPSO runs were appended.
#
# Excel's COM automation model expresses column widths in "number of
# characters" via Range.ColumnWidth / EntireColumn.ColumnWidth, which is
# exactly the same unit Excel itself uses when it regenerates the <col
# width="..."/> attribute on save, so we simply set that property on
# column A of the "Dados" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados")

$ws.Columns.Item(1).ColumnWidth = 16
